$wb = $excel.ActiveWorkbook

# --- Sheet "总计": insert a new row for the 2022-Q3 summary entry,
#     shifting the existing 2022-Q2 / 2021-Q2 rows down by one ---
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(3).Insert()

# Copy the formatting of row 2 onto the newly inserted row 3
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)  # xlPasteFormats

# Row 2 keeps the 2022-Q2 values but is relabeled 2022-Q3
$total.Range("B2").Value = "2022-Q3"

# Row 3 becomes the (old) 2022-Q2 summary row
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.01

# Row 4 (old row 3) is the 2021-Q2 summary row; only its index changes
$total.Range("A4").Value = 2

# --- Duplicate the "2022-Q2" sheet to create the new "2022-Q3" sheet,
#     inserting it immediately before "2022-Q2" ---
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# Update the fund-position figures on the new 2022-Q3 sheet.
# Columns D:G hold text-like numbers, so force text formatting before
# assigning them (otherwise the numeric-looking strings get auto-converted
# to real numbers), then drop back to the default style afterwards.
$q3.Range("D2:G3").NumberFormat = "@"

$q3.Range("D2").Value = "0.13"
$q3.Range("E2").Value = "92.87"
$q3.Range("F2").Value = "4.39"
$q3.Range("G2").Value = "0.0057"
$q3.Range("H2").Value = 4

$q3.Range("D3").Value = "0.08"
$q3.Range("E3").Value = "92.87"
$q3.Range("F3").Value = "4.39"
$q3.Range("G3").Value = "0.0035"
$q3.Range("H3").Value = 4

$q3.Range("D2:G3").Style = "Normal"

# Restore the original active/selected sheet (2021-Q2), since copying a
# sheet makes the copy the active one.
$wb.Worksheets.Item("2021-Q2").Activate()
